# "Add files via upload" — replaces the placeholder example-ticker row (row 2)
# with real stock tickers for each sector listed in row 1, and moves the
# active selection to G3.
#
# Column -> Sector (row 1, unchanged) -> new example ticker (row 2):
#   A  Basic Materials        -> dd     (DuPont, lowercase as typed)
#   B  Communication Services -> DD     (DuPont)
#   C  Consumer Cyclical      -> VZ     (Verizon)
#   D  Consumer Defensive     -> MCD    (McDonald's)
#   E  Energy                 -> XOM    (Exxon Mobil)
#   F  Financial Services     -> BRK-B  (Berkshire Hathaway)
#   G  Health Care            -> LLY    (Eli Lilly)
#   H  Industrials            -> LMT    (Lockheed Martin)
#   I  Real Estate            -> SPG    (Simon Property Group)
#   J  Technology             -> IBM
#   K  Utilities              -> SO     (Southern Company)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dd"
$ws.Range("B2").Value = "DD"
$ws.Range("C2").Value = "VZ"
$ws.Range("D2").Value = "MCD"
$ws.Range("E2").Value = "XOM"
$ws.Range("F2").Value = "BRK-B"
$ws.Range("G2").Value = "LLY"
$ws.Range("H2").Value = "LMT"
$ws.Range("I2").Value = "SPG"
$ws.Range("J2").Value = "IBM"
$ws.Range("K2").Value = "SO"

# Move the active selection to G3, matching the saved view state in the
# upload.
$ws.Range("G3").Select()
